$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @("NSAA", "sensorMagneticField", "dhc", 60),
    @("NSAA", "sensorMagneticField", "dhc", 60),
    @("NSAA", "jointAngle", "dhc", 60),
    @("NSAA", "AD", "dhc", 10),
    @("NSAA", "sensorMagneticField", "dhc", 60),
    @("NSAA", "sensorMagneticField", "overall", 60),
    @("NSAA", "sensorMagneticField", "acts", 60),
    @("NSAA", "jointAngle", "dhc", 60)
)

$startRow = 338
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}
